# Atualizacao de bases das ligas, do dia: 03-04-2024 as 22:09
# Several match rows had their data (all columns except the sequential
# index in column A) swapped with an adjacent row. Re-create that swap
# by exchanging the B:AC ranges of each affected row pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(129, 130),
    @(154, 156),
    @(157, 158),
    @(168, 169),
    @(210, 211)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AC${r1}")
    $range2 = $ws.Range("B${r2}:AC${r2}")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
